# Updates the crypto price/volume snapshot (and swaps the ONDO/Cosmos row
# order) to match the GitHub Actions refresh of Fri Jul 12 08:32:50 UTC 2024.
#
# Price cells in column D are stored as text in the workbook (e.g. "57.31",
# "0.998") even though they look numeric, so plain numeric-looking values
# would otherwise get auto-coerced into real numbers by Excel. For those we
# assign with a leading apostrophe (forces text entry) and then reset the
# cell style to "Normal" so the transient quote-prefix formatting introduced
# by the apostrophe doesn't linger as a style on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.139.63"
$ws.Range("E2").Value = "  -1.74%  "
$ws.Range("D3").Value = "3.072.29"
$ws.Range("E3").Value = "  -1.70%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'521.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.22%  "
$ws.Range("D6").Value = "'135.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.01%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "3.071.59"
$ws.Range("E8").Value = "  -1.71%  "
$ws.Range("E9").Value = "  +4.21%  "
$ws.Range("D10").Value = "'7.31"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.74%  "
$ws.Range("E11").Value = "  -2.24%  "
$ws.Range("E12").Value = "  +1.50%  "
$ws.Range("E13").Value = "  +1.46%  "
$ws.Range("D14").Value = "3.603.78"
$ws.Range("E14").Value = "  -1.52%  "
$ws.Range("D15").Value = "'25.22"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.87%  "
$ws.Range("D16").Value = "'0.0000160"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.82%  "
$ws.Range("D17").Value = "57.170.72"
$ws.Range("E17").Value = "  -1.78%  "
$ws.Range("D18").Value = "3.080.64"
$ws.Range("E18").Value = "  -1.24%  "
$ws.Range("D19").Value = "'5.86"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.23%  "
$ws.Range("D20").Value = "'12.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.97%  "
$ws.Range("D21").Value = "'7.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.27%  "
$ws.Range("D22").Value = "'349.79"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.05%  "
$ws.Range("D23").Value = "'0.998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("D24").Value = "'68.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.74%  "
$ws.Range("D25").Value = "'0.497"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.41%  "
$ws.Range("E26").Value = "  -2.38%  "
$ws.Range("D27").Value = "'1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.33%  "
$ws.Range("D28").Value = "0.0₃0865"
$ws.Range("E28").Value = "  -6.36%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").Value = "'7.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.77%  "
$ws.Range("E31").Value = "  -1.13%  "
$ws.Range("D32").Value = "'5.83"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.52%  "
$ws.Range("D33").Value = "'20.89"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.96%  "
$ws.Range("E34").Value = "  +2.80%  "
$ws.Range("E35").Value = "  +0.28%  "
$ws.Range("E36").Value = "  -5.48%  "
$ws.Range("D37").Value = "'5.98"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.41%  "
$ws.Range("D38").Value = "'25.46"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.37%  "
$ws.Range("D39").Value = "'1.22"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Value = "'0.0654"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.27%  "
$ws.Range("D41").Value = "'1.57"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.40%  "
$ws.Range("D42").Value = "'4.06"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.40%  "
$ws.Range("D43").Value = "'0.691"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.34%  "
$ws.Range("D44").Value = "2.399.32"
$ws.Range("E44").Value = "  +5.19%  "
$ws.Range("D45").Value = "'36.61"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("D46").Value = "'0.999"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.07%  "
$ws.Range("D47").Value = "3.113.70"
$ws.Range("E47").Value = "  -1.56%  "
$ws.Range("E48").Value = "  -0.70%  "
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").Value = "'5.96"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.46%  "
$ws.Range("B50").Value = "ONDO"
$ws.Range("C50").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D50").Value = "'0.942"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.00%  "
$ws.Range("D51").Value = "'19.55"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.34%  "
